$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13714
$ws1.Range("F4").Value = 13494
$ws1.Range("F6").Value = 802
$ws1.Range("F8").Value = 588
$ws1.Range("F20").Value = 423
$ws1.Range("F22").Value = 315
$ws1.Range("F23").Value = 215
$ws1.Range("F24").Value = 822
$ws1.Range("F26").Value = 1

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 158
$ws2.Range("F7").Value = 1444

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 13714
$ws4.Range("F5").Value = 13494
$ws4.Range("F7").Value = 802
$ws4.Range("F9").Value = 588
$ws4.Range("F27").Value = 423
$ws4.Range("F29").Value = 315
$ws4.Range("F30").Value = 215
$ws4.Range("F31").Value = 822
$ws4.Range("F32").Value = 158
$ws4.Range("F33").Value = 1444
$ws4.Range("F39").Value = 1
